$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1883.3334
$ws.Range("J112").Value = 2135.6428
$ws.Range("L112").Value = 6406.928400000001
$ws.Range("N112").Value = -8622.928400000001
$ws.Range("H129").Value = 836774.5
$ws.Range("I129").Value = 472.66666
$ws.Range("J129").Value = 1079571.9
$ws.Range("K129").Value = 1417.99998
$ws.Range("L129").Value = 3238715.7
$ws.Range("M129").Value = 3582.00002
$ws.Range("N129").Value = -3248715.7
$ws.Range("H132").Value = 2423324.8
$ws.Range("I132").Value = 2647521.2
$ws.Range("K132").Value = 7942563.600000001
$ws.Range("M132").Value = -7940033.600000001
$ws.Range("H137").Value = 20409694
$ws.Range("I137").Value = 1162.0294
$ws.Range("K137").Value = 3486.0882
$ws.Range("M137").Value = -936.0881999999997
$ws.Range("H138").Value = 3606.0476
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3606.0476
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 10818.1428
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -21098.1428
$ws.Range("H141").Value = 3201.7273
$ws.Range("I141").Value = 2646.9
$ws.Range("J141").Value = 8750
$ws.Range("K141").Value = 7940.700000000001
$ws.Range("L141").Value = 26250
$ws.Range("M141").Value = -2760.700000000001
$ws.Range("N141").Value = -36610

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1113.3636
$ws.Range("I2").Value = 576.5
$ws.Range("J2").Value = 1420.1428
$ws.Range("K2").Value = 576.5
$ws.Range("L2").Value = 1420.1428
$ws.Range("M2").Value = -463.5
$ws.Range("N2").Value = -1646.1428
$ws.Range("H32").Value = 42035.535
$ws.Range("I32").Value = 46134.785
$ws.Range("J32").Value = 33206.383
$ws.Range("K32").Value = 46134.785
$ws.Range("L32").Value = 33206.383
$ws.Range("M32").Value = -45847.785
$ws.Range("N32").Value = -33780.383
$ws.Range("H45").Value = 825.53845
$ws.Range("I45").Value = 795.9091
$ws.Range("J45").Value = 988.5
$ws.Range("K45").Value = 795.9091
$ws.Range("L45").Value = 988.5
$ws.Range("M45").Value = -418.9091
$ws.Range("N45").Value = -1742.5
$ws.Range("H54").Value = 28000
$ws.Range("J54").Value = 28000
$ws.Range("L54").Value = 28000
$ws.Range("N54").Value = -29538
$ws.Range("H74").Value = 874.3134
$ws.Range("I74").Value = 724.8570999999999
$ws.Range("J74").Value = 1635.1818
$ws.Range("K74").Value = 724.8570999999999
$ws.Range("L74").Value = 1635.1818
$ws.Range("M74").Value = 149.1429000000001
$ws.Range("N74").Value = -3383.1818
$ws.Range("H77").Value = 874.3134
$ws.Range("I77").Value = 724.8570999999999
$ws.Range("J77").Value = 1635.1818
$ws.Range("K77").Value = 3624.2855
$ws.Range("L77").Value = 8175.909000000001
$ws.Range("M77").Value = 743.7145
$ws.Range("N77").Value = -16911.909
$ws.Range("H80").Value = 27000
$ws.Range("J80").Value = 27000
$ws.Range("L80").Value = 27000
$ws.Range("N80").Value = -28996
$ws.Range("H83").Value = 27000
$ws.Range("J83").Value = 27000
$ws.Range("L83").Value = 81000
$ws.Range("N83").Value = -90984
$ws.Range("H116").Value = 1113.3636
$ws.Range("I116").Value = 576.5
$ws.Range("J116").Value = 1420.1428
$ws.Range("K116").Value = 576.5
$ws.Range("L116").Value = 1420.1428
$ws.Range("M116").Value = 1717.5
$ws.Range("N116").Value = -6008.1428
$ws.Range("H132").Value = 6706.8477
$ws.Range("I132").Value = 7948.8184
$ws.Range("J132").Value = 3554.1538
$ws.Range("K132").Value = 23846.4552
$ws.Range("L132").Value = 10662.4614
$ws.Range("M132").Value = -21316.4552
$ws.Range("N132").Value = -15722.4614

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1113.3636
$ws.Range("I3").Value = 576.5
$ws.Range("J3").Value = 1420.1428
$ws.Range("K3").Value = 576.5
$ws.Range("L3").Value = 1420.1428
$ws.Range("M3").Value = -462.5
$ws.Range("N3").Value = -1648.1428
$ws.Range("H22").Value = 286.09525
$ws.Range("I22").Value = 280.4
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 280.4
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -107.4
$ws.Range("N22").Value = -746
$ws.Range("H134").Value = 3485.2603
$ws.Range("I134").Value = 3635.8518
$ws.Range("J134").Value = 3057.2632
$ws.Range("K134").Value = 10907.5554
$ws.Range("L134").Value = 9171.7896
$ws.Range("M134").Value = -8372.555399999999
$ws.Range("N134").Value = -14241.7896

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 905.375
$ws.Range("I16").Value = 802.3333
$ws.Range("J16").Value = 1037.8572
$ws.Range("K16").Value = 802.3333
$ws.Range("L16").Value = 1037.8572
$ws.Range("M16").Value = -515.3333
$ws.Range("N16").Value = -1611.8572
$ws.Range("H31").Value = 16953520
$ws.Range("I31").Value = 2063.2812
$ws.Range("K31").Value = 2063.2812
$ws.Range("M31").Value = -1768.2812
$ws.Range("H34").Value = 16953520
$ws.Range("I34").Value = 2063.2812
$ws.Range("K34").Value = 2063.2812
$ws.Range("M34").Value = -1861.2812
$ws.Range("H113").Value = 905.375
$ws.Range("I113").Value = 802.3333
$ws.Range("J113").Value = 1037.8572
$ws.Range("K113").Value = 802.3333
$ws.Range("L113").Value = 1037.8572
$ws.Range("M113").Value = 1367.6667
$ws.Range("N113").Value = -5377.8572
$ws.Range("H132").Value = 5001834
$ws.Range("I132").Value = 1532
$ws.Range("J132").Value = 25003042
$ws.Range("K132").Value = 4596
$ws.Range("L132").Value = 75009126
$ws.Range("M132").Value = -2066
$ws.Range("N132").Value = -75014186
$ws.Range("H134").Value = 3596.0557
$ws.Range("I134").Value = 3704.2354
$ws.Range("J134").Value = 1757
$ws.Range("K134").Value = 11112.7062
$ws.Range("L134").Value = 5271
$ws.Range("M134").Value = -8577.706200000001
$ws.Range("N134").Value = -10341

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 1686.25
$ws.Range("H131").Value = 1503870.8
$ws.Range("I131").Value = 8650
$ws.Range("J131").Value = 1746844.1
$ws.Range("K131").Value = 25950
$ws.Range("L131").Value = 5240532.300000001
$ws.Range("M131").Value = -20910
$ws.Range("N131").Value = -5250612.300000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 42000
$ws.Range("J69").Value = 42000
$ws.Range("L69").Value = 42000
$ws.Range("N69").Value = -43498
$ws.Range("H72").Value = 42000
$ws.Range("J72").Value = 42000
$ws.Range("L72").Value = 126000
$ws.Range("N72").Value = -133488

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 841
$ws.Range("H27").Value = 841
$ws.Range("H46").Value = 1453.1
$ws.Range("J46").Value = 1525.6666
$ws.Range("L46").Value = 1525.6666
$ws.Range("N46").Value = -1901.6666
$ws.Range("H55").Value = 465.75
$ws.Range("I55").Value = 500.25
$ws.Range("J55").Value = 431.25
$ws.Range("K55").Value = 500.25
$ws.Range("L55").Value = 431.25
$ws.Range("M55").Value = -327.25
$ws.Range("N55").Value = -777.25
$ws.Range("H132").Value = 16915.54
$ws.Range("I132").Value = 28072
$ws.Range("J132").Value = 3899.6667
$ws.Range("K132").Value = 84216
$ws.Range("L132").Value = 11699.0001
$ws.Range("M132").Value = -81686
$ws.Range("N132").Value = -16759.0001
$ws.Range("H136").Value = 9377.875
$ws.Range("I136").Value = 9803.066000000001
$ws.Range("K136").Value = 29409.198
$ws.Range("M136").Value = -26859.198

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H62").Value = 2361
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2361
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2361
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -3609
$ws.Range("H65").Value = 2361
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2361
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 11805
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -18045
